$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

$oldText = " When using the data, please ensure proper attribution to original paper present in the same github repository (Cupolillo et al. 2024 Movement pattern of Trinomys dimidiatus)."
$newText = " When using the data, please ensure proper attribution to original paper"

$find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
